$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.032.38"
$ws.Range("E2").Value = "  +1.45%  "

# Row 3
$ws.Range("D3").Value = "3.166.96"
$ws.Range("E3").Value = "  +3.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.91"
$ws.Range("E5").Value = "  +2.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.43"
$ws.Range("E6").Value = "  +4.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "3.165.46"
$ws.Range("E8").Value = "  +3.37%  "

# Row 9
$ws.Range("E9").Value = "  +2.64%  "

# Row 10
$ws.Range("E10").Value = "  +3.84%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.21"
$ws.Range("E11").Value = "  +1.83%  "

# Row 12
$ws.Range("E12").Value = "  +4.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +17.74%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.03"
$ws.Range("E14").Value = "  +6.84%  "

# Row 15
$ws.Range("D15").Value = "3.681.35"
$ws.Range("E15").Value = "  +3.43%  "

# Row 16
$ws.Range("D16").Value = "65.118.13"
$ws.Range("E16").Value = "  +1.64%  "

# Row 17
$ws.Range("D17").Value = "3.173.58"
$ws.Range("E17").Value = "  +3.65%  "

# Row 18
$ws.Range("E18").Value = "  +6.06%  "

# Row 19
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.37"
$ws.Range("E20").Value = "  +6.44%  "

# Row 21
$ws.Range("E21").Value = "  +6.44%  "

# Row 22
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.68"
$ws.Range("E22").Value = "  +8.19%  "

# Row 23
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("E23").Value = "  +7.45%  "

# Row 24
$ws.Range("E24").Value = "  +2.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.74"
$ws.Range("E25").Value = "  +3.13%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("E27").Value = "  +14.50%  "

# Row 28
$ws.Range("E28").Value = "  +3.47%  "

# Row 29
$ws.Range("E29").Value = "  +7.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.02"
$ws.Range("E30").Value = "  +6.47%  "

# Row 31
$ws.Range("E31").Value = "  +14.86%  "

# Row 32
$ws.Range("E32").Value = "  +7.15%  "

# Row 33
$ws.Range("E33").Value = "  +0.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.30"
$ws.Range("E34").Value = "  +10.58%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.67"
$ws.Range("E35").Value = "  +6.75%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.62"
$ws.Range("E36").Value = "  +1.30%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "480.02"
$ws.Range("E37").Value = "  +7.76%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0886"
$ws.Range("E38").Value = "  +9.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  +7.90%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0421"
$ws.Range("E40").Value = "  +2.79%  "

# Row 41
$ws.Range("D41").Value = "3.123.98"
$ws.Range("E41").Value = "  +4.41%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.61"
$ws.Range("E42").Value = "  +4.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.120"
$ws.Range("E43").Value = "  +4.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  +16.70%  "

# Row 45
$ws.Range("E45").Value = "  +10.43%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.22"
$ws.Range("E46").Value = "  +4.28%  "

# Row 47
$ws.Range("D47").Value = "0.0₃0593"
$ws.Range("E47").Value = "  +13.68%  "

# Row 50
$ws.Range("E50").Value = "  +10.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.54"
$ws.Range("E51").Value = "  +3.04%  "
